$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two more year columns (2022, 2023) to the table, following the
# --- same pattern already used by the 2012..2021 columns (D..M). Column M
# --- (2021) is the right-most existing column, so its formatting (font,
# --- border, number format) is the template for the new N/O columns.

# Row 4 (year headers): copy M4's format onto N4:O4, then fill in the years.
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4:O4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = 2022
$ws.Range("O4").Value = 2023

# Row 5 (data values): copy M5's format onto N5:O5, then repeat the same
# 6.53 value already used for 2019-2021 (K5:M5).
$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5:O5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").Value = 6.53
$ws.Range("O5").Value = 6.53

$excel.CutCopyMode = $false

# --- Row-height touch-ups that came along with the new columns.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 17.25
